$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(6, 8).Value = 754.63635
$ws.Cells.Item(6, 9).Value = 558.44446
$ws.Cells.Item(6, 11).Value = 1675.33338
$ws.Cells.Item(6, 13).Value = -1563.33338
$ws.Cells.Item(28, 8).Value = 2113.5
$ws.Cells.Item(28, 9).Value = 1884.8
$ws.Cells.Item(28, 10).Value = 2494.6667
$ws.Cells.Item(28, 11).Value = 1884.8
$ws.Cells.Item(28, 12).Value = 2494.6667
$ws.Cells.Item(28, 13).Value = -1399.8
$ws.Cells.Item(28, 14).Value = -3464.6667
$ws.Cells.Item(31, 8).Value = 3400
$ws.Cells.Item(31, 9).Value = 6000
$ws.Cells.Item(31, 11).Value = 18000
$ws.Cells.Item(31, 13).Value = -17770
$ws.Cells.Item(40, 8).Value = 2574.875
$ws.Cells.Item(40, 9).Value = 1766.6666
$ws.Cells.Item(40, 11).Value = 1766.6666
$ws.Cells.Item(40, 13).Value = -1591.6666
$ws.Cells.Item(76, 8).Value = 3528
$ws.Cells.Item(76, 9).Value = 3528
$ws.Cells.Item(76, 11).Value = 3528
$ws.Cells.Item(76, 13).Value = -3213
$ws.Cells.Item(79, 8).Value = 3528
$ws.Cells.Item(79, 9).Value = 3528
$ws.Cells.Item(79, 11).Value = 3528
$ws.Cells.Item(79, 13).Value = -2436
$ws.Cells.Item(116, 8).Value = 1131171
$ws.Cells.Item(116, 9).Value = 1562985.4
$ws.Cells.Item(116, 10).Value = 8453.6
$ws.Cells.Item(116, 11).Value = 1562985.4
$ws.Cells.Item(116, 12).Value = 8453.6
$ws.Cells.Item(116, 13).Value = -1559543.4
$ws.Cells.Item(116, 14).Value = -15337.6
$ws.Cells.Item(121, 8).Value = 3600
$ws.Cells.Item(121, 10).Value = 3600
$ws.Cells.Item(121, 12).Value = 10800
$ws.Cells.Item(121, 14).Value = -14294
$ws.Cells.Item(132, 8).Value = 3502881.5
$ws.Cells.Item(132, 9).Value = 3779330
$ws.Cells.Item(132, 11).Value = 11337990
$ws.Cells.Item(132, 13).Value = -11335460

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(3, 8).Value = 1000
$ws.Cells.Item(3, 9).Value = 1000
$ws.Cells.Item(3, 11).Value = 1000
$ws.Cells.Item(3, 13).Value = -885
$ws.Cells.Item(45, 8).Value = 2367.0417
$ws.Cells.Item(45, 9).Value = 1186.1666
$ws.Cells.Item(45, 10).Value = 3547.9167
$ws.Cells.Item(45, 11).Value = 1186.1666
$ws.Cells.Item(45, 12).Value = 3547.9167
$ws.Cells.Item(45, 13).Value = -809.1666
$ws.Cells.Item(45, 14).Value = -4301.9167
$ws.Cells.Item(61, 8).Value = 10072
$ws.Cells.Item(61, 9).Value = 1137.7142
$ws.Cells.Item(61, 10).Value = 22580
$ws.Cells.Item(61, 11).Value = 1137.7142
$ws.Cells.Item(61, 12).Value = 22580
$ws.Cells.Item(61, 13).Value = -925.7141999999999
$ws.Cells.Item(61, 14).Value = -23004
$ws.Cells.Item(74, 8).Value = 612642.3
$ws.Cells.Item(74, 9).Value = 6000012
$ws.Cells.Item(74, 11).Value = 6000012
$ws.Cells.Item(74, 13).Value = -5999138
$ws.Cells.Item(77, 8).Value = 612642.3
$ws.Cells.Item(77, 9).Value = 6000012
$ws.Cells.Item(77, 11).Value = 30000060
$ws.Cells.Item(77, 13).Value = -29995692
$ws.Cells.Item(97, 8).Value = 809.7174
$ws.Cells.Item(97, 9).Value = 677.381
$ws.Cells.Item(97, 11).Value = 677.381
$ws.Cells.Item(97, 13).Value = -181.381
$ws.Cells.Item(110, 8).Value = 4880.8276
$ws.Cells.Item(110, 9).Value = 4653.423
$ws.Cells.Item(110, 11).Value = 4653.423
$ws.Cells.Item(110, 13).Value = -2608.423
$ws.Cells.Item(122, 8).Value = 1606.48
$ws.Cells.Item(122, 9).Value = 1388.238
$ws.Cells.Item(122, 11).Value = 4164.714
$ws.Cells.Item(122, 13).Value = -1714.714
$ws.Cells.Item(132, 8).Value = 1511.8214
$ws.Cells.Item(132, 9).Value = 752.65
$ws.Cells.Item(132, 11).Value = 2257.95
$ws.Cells.Item(132, 13).Value = 272.0500000000002
$ws.Cells.Item(136, 8).Value = 10072
$ws.Cells.Item(136, 9).Value = 1137.7142
$ws.Cells.Item(136, 10).Value = 22580
$ws.Cells.Item(136, 11).Value = 3413.1426
$ws.Cells.Item(136, 12).Value = 67740
$ws.Cells.Item(136, 13).Value = -863.1425999999997
$ws.Cells.Item(136, 14).Value = -72840

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 1248.7084
$ws.Cells.Item(20, 9).Value = 934.4
$ws.Cells.Item(20, 10).Value = 1772.5555
$ws.Cells.Item(20, 11).Value = 934.4
$ws.Cells.Item(20, 12).Value = 1772.5555
$ws.Cells.Item(20, 13).Value = -687.4
$ws.Cells.Item(20, 14).Value = -2266.5555
$ws.Cells.Item(82, 8).Value = 20367.889
$ws.Cells.Item(82, 9).Value = 13051.833
$ws.Cells.Item(82, 11).Value = 13051.833
$ws.Cells.Item(82, 13).Value = -12668.833
$ws.Cells.Item(85, 8).Value = 20367.889
$ws.Cells.Item(85, 9).Value = 13051.833
$ws.Cells.Item(85, 11).Value = 13051.833
$ws.Cells.Item(85, 13).Value = -11725.833
$ws.Cells.Item(99, 8).Value = 1499.8334
$ws.Cells.Item(99, 9).Value = 1499.8334
$ws.Cells.Item(99, 11).Value = 1499.8334
$ws.Cells.Item(99, 13).Value = -1.833399999999983
$ws.Cells.Item(105, 8).Value = 2397.6365
$ws.Cells.Item(105, 9).Value = 2391.4707
$ws.Cells.Item(105, 10).Value = 2418.6
$ws.Cells.Item(105, 11).Value = 2391.4707
$ws.Cells.Item(105, 12).Value = 2418.6
$ws.Cells.Item(105, 13).Value = -644.4706999999999
$ws.Cells.Item(105, 14).Value = -5912.6

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(10, 8).Value = 10666
$ws.Cells.Item(10, 9).Value = 0
$ws.Cells.Item(10, 10).Value = 10666
$ws.Cells.Item(10, 11).Value = 0
$ws.Cells.Item(10, 12).ClearContents()
$ws.Cells.Item(10, 13).Value = 10666
$ws.Cells.Item(10, 14).Value = -10944
$ws.Cells.Item(22, 8).Value = 665.6667
$ws.Cells.Item(22, 10).Value = 747.1429000000001
$ws.Cells.Item(22, 12).Value = 747.1429000000001
$ws.Cells.Item(22, 14).Value = -1447.1429
$ws.Cells.Item(86, 8).Value = 80209.44500000001
$ws.Cells.Item(86, 9).Value = 111147.836
$ws.Cells.Item(86, 10).Value = 18332.666
$ws.Cells.Item(86, 11).Value = 111147.836
$ws.Cells.Item(86, 12).Value = 18332.666
$ws.Cells.Item(86, 13).Value = -110024.836
$ws.Cells.Item(86, 14).Value = -20578.666
$ws.Cells.Item(89, 8).Value = 80209.44500000001
$ws.Cells.Item(89, 9).Value = 111147.836
$ws.Cells.Item(89, 10).Value = 18332.666
$ws.Cells.Item(89, 11).Value = 555739.1799999999
$ws.Cells.Item(89, 12).Value = 91663.33
$ws.Cells.Item(89, 13).Value = -550123.1799999999
$ws.Cells.Item(89, 14).Value = -102895.33

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(11, 8).Value = 2960927.8
$ws.Cells.Item(11, 9).Value = 3806742.8
$ws.Cells.Item(11, 10).Value = 575
$ws.Cells.Item(11, 11).Value = 11420228.4
$ws.Cells.Item(11, 12).Value = 1725
$ws.Cells.Item(11, 13).Value = -11420088.4
$ws.Cells.Item(11, 14).Value = -2005
$ws.Cells.Item(107, 8).Value = 1417.8182
$ws.Cells.Item(107, 9).Value = 1844.25
$ws.Cells.Item(107, 10).Value = 906.1
$ws.Cells.Item(107, 11).Value = 5532.75
$ws.Cells.Item(107, 12).Value = 2718.3
$ws.Cells.Item(107, 13).Value = -3612.75
$ws.Cells.Item(107, 14).Value = -6558.3
$ws.Cells.Item(123, 8).Value = 4450
$ws.Cells.Item(123, 9).Value = 0
$ws.Cells.Item(123, 11).Value = 0
$ws.Cells.Item(123, 13).ClearContents()

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 5754
$ws.Cells.Item(70, 9).Value = 7962.6665
$ws.Cells.Item(70, 10).Value = 4649.6665
$ws.Cells.Item(70, 11).Value = 7962.6665
$ws.Cells.Item(70, 12).Value = 4649.6665
$ws.Cells.Item(70, 13).Value = -7692.6665
$ws.Cells.Item(70, 14).Value = -5189.6665
$ws.Cells.Item(73, 8).Value = 5754
$ws.Cells.Item(73, 9).Value = 7962.6665
$ws.Cells.Item(73, 10).Value = 4649.6665
$ws.Cells.Item(73, 11).Value = 7962.6665
$ws.Cells.Item(73, 12).Value = 4649.6665
$ws.Cells.Item(73, 13).Value = -7026.6665
$ws.Cells.Item(73, 14).Value = -6521.6665
$ws.Cells.Item(80, 8).Value = 10391.538
$ws.Cells.Item(80, 9).Value = 4699.1665
$ws.Cells.Item(80, 11).Value = 4699.1665
$ws.Cells.Item(80, 13).Value = -3701.1665
$ws.Cells.Item(83, 8).Value = 10391.538
$ws.Cells.Item(83, 9).Value = 4699.1665
$ws.Cells.Item(83, 11).Value = 23495.8325
$ws.Cells.Item(83, 13).Value = -18503.8325
$ws.Cells.Item(132, 8).Value = 1739.1
$ws.Cells.Item(132, 9).Value = 1567.4736
$ws.Cells.Item(132, 11).Value = 4702.4208
$ws.Cells.Item(132, 13).Value = -2172.4208

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(46, 8).Value = 2756.52
$ws.Cells.Item(46, 9).Value = 1110.7142
$ws.Cells.Item(46, 10).Value = 4851.1816
$ws.Cells.Item(46, 11).Value = 1110.7142
$ws.Cells.Item(46, 12).Value = 4851.1816
$ws.Cells.Item(46, 13).Value = -922.7141999999999
$ws.Cells.Item(46, 14).Value = -5227.1816
$ws.Cells.Item(61, 8).Value = 2927
$ws.Cells.Item(61, 9).Value = 2927
$ws.Cells.Item(61, 10).Value = 0
$ws.Cells.Item(61, 11).Value = 2927
$ws.Cells.Item(61, 12).Value = 0
$ws.Cells.Item(61, 13).ClearContents()
$ws.Cells.Item(61, 14).Value = -2725
$ws.Cells.Item(68, 8).Value = 2829.5264
$ws.Cells.Item(68, 9).Value = 2455.5
$ws.Cells.Item(68, 10).Value = 3470.7144
$ws.Cells.Item(68, 11).Value = 2455.5
$ws.Cells.Item(68, 12).Value = 3470.7144
$ws.Cells.Item(68, 13).Value = -1706.5
$ws.Cells.Item(68, 14).Value = -4968.7144
$ws.Cells.Item(71, 8).Value = 2829.5264
$ws.Cells.Item(71, 9).Value = 2455.5
$ws.Cells.Item(71, 10).Value = 3470.7144
$ws.Cells.Item(71, 11).Value = 12277.5
$ws.Cells.Item(71, 12).Value = 17353.572
$ws.Cells.Item(71, 13).Value = -8533.5
$ws.Cells.Item(71, 14).Value = -24841.572
$ws.Cells.Item(113, 8).Value = 2927
$ws.Cells.Item(113, 9).Value = 2927
$ws.Cells.Item(113, 10).Value = 0
$ws.Cells.Item(113, 11).Value = 2927
$ws.Cells.Item(113, 12).Value = 0
$ws.Cells.Item(113, 13).ClearContents()
$ws.Cells.Item(113, 14).Value = -757
$ws.Cells.Item(136, 8).Value = 2921.2856
$ws.Cells.Item(136, 9).Value = 2885.4443
$ws.Cells.Item(136, 10).Value = 2985.8
$ws.Cells.Item(136, 11).Value = 8656.332900000001
$ws.Cells.Item(136, 12).Value = 8957.400000000001
$ws.Cells.Item(136, 13).Value = -6106.332900000001
$ws.Cells.Item(136, 14).Value = -14057.4

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(37, 8).Value = 0
$ws.Cells.Item(37, 10).Value = 0
$ws.Cells.Item(37, 12).ClearContents()
$ws.Cells.Item(37, 14).Value = 0
$ws.Cells.Item(62, 8).Value = 3662.75
$ws.Cells.Item(62, 9).Value = 3260
$ws.Cells.Item(62, 11).Value = 3260
$ws.Cells.Item(62, 13).Value = -2636
$ws.Cells.Item(65, 8).Value = 3662.75
$ws.Cells.Item(65, 9).Value = 3260
$ws.Cells.Item(65, 11).Value = 16300
$ws.Cells.Item(65, 13).Value = -13180
$ws.Cells.Item(122, 8).Value = 24783.771
$ws.Cells.Item(122, 9).Value = 27056.71
$ws.Cells.Item(122, 11).Value = 81170.13
$ws.Cells.Item(122, 13).Value = -78720.13
$ws.Cells.Item(132, 8).Value = 35134.81
$ws.Cells.Item(132, 9).Value = 42875.117
$ws.Cells.Item(132, 10).Value = 2238.5
$ws.Cells.Item(132, 11).Value = 128625.351
$ws.Cells.Item(132, 12).Value = 6715.5
$ws.Cells.Item(132, 13).Value = -126095.351
$ws.Cells.Item(132, 14).Value = -11775.5
